$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A84").Value = "2025-04-29 13:46:35"
$ws.Range("B84").Value = 269
